$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Activate()
Write-Host $wb.Windows.Count
$win = $wb.Windows.Item(1)
$win.ScrollRow = 10
$win.ScrollColumn = 1
Write-Host $win.ScrollRow
